$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SBUX")

# Row 15 - Gross Margin
$ws.Range("D15").Value = 0.2752
$ws.Range("E15").Value = 0.4156
$ws.Range("F15").Value = 0.5595
$ws.Range("G15").Value = 0.6816

# Row 19 - Free Cash Flow Margin
$ws.Range("D19").Value = 0.0444
$ws.Range("E19").Value = 0.007
$ws.Range("F19").Value = 0.0503
$ws.Range("G19").Value = 0.1014

# Row 27 - EBITDA Margin
$ws.Range("D27").Value = 0.1569
$ws.Range("E27").Value = 0.1589
$ws.Range("F27").Value = 0.1994
$ws.Range("G27").Value = 0.2132

# Row 28 - Operating Cash Flow Margin
$ws.Range("D28").Value = 0.1075
$ws.Range("E28").Value = 0.0762
$ws.Range("F28").Value = 0.1148
$ws.Range("G28").Value = 0.167
